$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 213 (shifts existing rows 213:338 down to 214:339,
# extending the used range to A1:R339) and fill it with the new weekly
# price record.
$ws.Rows("213:213").Insert()

$ws.Range("A213").Value = 5
$ws.Range("B213").Value = 'Macroferia Regional de Talca'
$ws.Range("C213").Value = 'Maule'
$ws.Range("D213").Value = 44879
$ws.Range("E213").Value = 7
$ws.Range("F213").Value = 100112009
$ws.Range("G213").Value = 'Acelga'
$ws.Range("H213").Value = 'Sin especificar'
$ws.Range("I213").Value = 'Primera'
$ws.Range("J213").Value = 500
$ws.Range("K213").Value = 2500
$ws.Range("L213").Value = 2500
$ws.Range("M213").Value = 2500
$ws.Range("N213").Value = '$/docena de atados (4 kilos)'
$ws.Range("O213").Value = 'Región del Maule'
$ws.Range("P213").Value = 625
$ws.Range("Q213").Value = 4
$ws.Range("R213").Value = 'Hortaliza'
